# fix: [#64663] Fix template
# Insert a new "Generalforsamlingsdato" column at the front of the header
# row and append a new "Nettoudbytte" column at the end, shifting the
# existing headers one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Generalforsamlingsdato"
$ws.Range("B1").Value = "Udbetalingsdato"
$ws.Range("C1").Value = "Identifikation"
$ws.Range("D1").Value = "Navn"
$ws.Range("E1").Value = "C/O"
$ws.Range("F1").Value = "Adresse"
$ws.Range("G1").Value = "Postnr."
$ws.Range("H1").Value = "Land"
$ws.Range("I1").Value = "Bruttoudbytte"
$ws.Range("J1").Value = "Nettoudbytte"

# New header cells should carry the same bold header styling as the rest
# of the row.
$ws.Range("I1:J1").Font.Bold = $true

$ws.Range("B2").Select() | Out-Null
